# Applies the "Abw. [σ]" header rename, the new "0.00" number formats,
# the refreshed border lines, the widened C:F / H:J columns, the
# landscape page orientation and the moved selection described in the
# commit "Im Jupyter Skript für 256 beta hinzugefügt".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text: "Abw." -> "Abw. [σ]" (both Kalpha and Kbeta blocks) ---
# K3 first, as a plain string (mirrors the plain shared-string entry used
# by the workbook for the second occurrence of the label).
$k3 = $ws.Range("K3")
$k3.Value = "Abw. [" + [char]0x03C3 + "]"

# G3 second, then re-flag the sigma glyph with an explicit font so it is
# stored as its own run (mirrors the mixed-font shared-string entry used
# for the first occurrence of the label).
$g3 = $ws.Range("G3")
$g3.Value = "Abw. [" + [char]0x03C3 + "]"
$g3.Characters(7, 1).Font.Name = "Calibri"

# --- 2. Number formats: exp./Delta/Abw. columns now show 2 decimals ---
$ws.Range("D4:D11").NumberFormat = "0.00"
$ws.Range("E4:E11").NumberFormat = "0.00"
$ws.Range("H4:H11").NumberFormat = "General"
$ws.Range("K4:K11").NumberFormat = "0.00"
$ws.Range("G4").NumberFormat = "0.00"
$ws.Range("G6:G9").NumberFormat = "0.00"
$ws.Range("G11").NumberFormat = "0.00"

# --- 3. Border touch-ups ---
# B2:C2 gain a thin bottom rule (matching the header block beneath them).
$ws.Range("B2:C2").Borders.Item(9).LineStyle = 1
$ws.Range("B2:C2").Borders.Item(9).Weight = 2

# B3 gains a medium right + bottom rule, closing off the row/column header box.
$ws.Range("B3").Borders.Item(10).LineStyle = 1
$ws.Range("B3").Borders.Item(10).Weight = -4138
$ws.Range("B3").Borders.Item(9).LineStyle = 1
$ws.Range("B3").Borders.Item(9).Weight = -4138

# C4:C11 right edge thickens from thin to medium ...
$ws.Range("C4:C11").Borders.Item(10).LineStyle = 1
$ws.Range("C4:C11").Borders.Item(10).Weight = -4138

# ... while D4:D11 loses its (now redundant) left edge.
$ws.Range("D4:D11").Borders.Item(7).LineStyle = -4142

# --- 4. Column widths: give the new narrower data columns a fixed width ---
$ws.Range("C1:F1").ColumnWidth = 6.3
$ws.Range("H1:J1").ColumnWidth = 6.3

# --- 5. Selection cursor moved ---
$ws.Range("E18").Select()

# --- 6. Page orientation: portrait -> landscape ---
$ws.PageSetup.Orientation = 2
